# "Updating the spawn rates" -------------------------------------------
# Populates the per-location "Frequency" matrix (columns Q:AA) on Sheet2
# with the new spawn-rate figures, fixes the General "Blockability" /
# "Escapability" figures that were mistakenly stored as fractions
# (0.3/0.6/0.4/0.7) instead of whole numbers, hides two rows that are no
# longer tracked (Feral cats / Pigeon), marks the still-unconfirmed
# Pigeon/Doves frequencies with "?", and appends COUNT/SUM roll-up rows
# at the bottom of the frequency table.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: new / updated frequency figures by location --------------
$ws2.Range("R3").Value  = 10
$ws2.Range("S3").Value  = 10
$ws2.Range("W3").Value  = 10
$ws2.Range("Y3").Value  = 2

$ws2.Range("T4").Value  = 3

$ws2.Range("U5").Value  = 5

$ws2.Range("S6").Value  = 10
$ws2.Range("V6").Value  = 8

$ws2.Range("Y7").Value  = 10

# General stats that were stored as fractions - correct to whole numbers
$ws2.Range("N8").Value  = 3
$ws2.Range("O8").Value  = 6

# Feral cats no longer tracked
$ws2.Rows.Item(10).Hidden = $true

$ws2.Range("R11").Value  = 2
$ws2.Range("S11").Value  = 10
$ws2.Range("U11").Value  = 10
$ws2.Range("AA11").Value = 10

$ws2.Range("S12").Value  = 8
$ws2.Range("V12").Value  = 6
$ws2.Range("W12").Value  = 3
$ws2.Range("X12").Value  = 8

$ws2.Range("R13").Value  = 1
$ws2.Range("S13").Value  = 10
$ws2.Range("V13").Value  = 8
$ws2.Range("AA13").Value = 8

$ws2.Range("S14").Value  = 10
$ws2.Range("W14").Value  = 3
$ws2.Range("Y14").Value  = 3

$ws2.Range("Q16").Value  = 8
$ws2.Range("R16").Value  = 8

$ws2.Range("R17").Value  = 10
$ws2.Range("S17").Value  = 8
$ws2.Range("V17").Value  = 10
$ws2.Range("Y17").Value  = 10
$ws2.Range("Z17").Value  = 10
$ws2.Range("AA17").Value = 8

$ws2.Range("AA18").Value = 10

$ws2.Range("T19").Value  = 1

$ws2.Range("Q20").Value  = 4
$ws2.Range("R20").Value  = 6

$ws2.Range("Q21").Value  = 10
$ws2.Range("R21").Value  = 10
$ws2.Range("U21").Value  = 3
$ws2.Range("AA21").Value = 6

# General stats that were stored as fractions - correct to whole numbers
$ws2.Range("N22").Value  = 4
$ws2.Range("O22").Value  = 7

# Pigeon no longer tracked; Pigeon & Doves frequencies unconfirmed ("?")
$ws2.Rows.Item(24).Hidden = $true
$ws2.Range("Q24:AA24").Value = "?"
$ws2.Range("Q25").Value = "?"

$ws2.Range("R26").Value  = 10
$ws2.Range("Y26").Value  = 10
$ws2.Range("Z26").Value  = 10

$ws2.Range("S27").Value  = 5
$ws2.Range("V27").Value  = 10
$ws2.Range("Y27").Value  = 10
$ws2.Range("Z27").Value  = 10

$ws2.Range("S28").Value  = 10
$ws2.Range("V28").Value  = 5
$ws2.Range("Z28").Value  = 5

$ws2.Range("S29").Value  = 10
$ws2.Range("V29").Value  = 5
$ws2.Range("Y29").Value  = 5
$ws2.Range("Z29").Value  = 8
$ws2.Range("AA29").Value = 3

$ws2.Range("S30").Value  = 3
$ws2.Range("W30").Value  = 3
$ws2.Range("X30").Value  = 3

$ws2.Range("T31").Value  = 1

$ws2.Range("U32").Value  = 10
$ws2.Range("AA32").Value = 5

$ws2.Range("S33").Value  = 3
$ws2.Range("Z33").Value  = 3
$ws2.Range("AA33").Value = 3

$ws2.Range("V34").Value  = 10
$ws2.Range("Z34").Value  = 10

$ws2.Range("X35").Value  = 10

$ws2.Range("N36").Value  = 5
$ws2.Range("W36").Value  = 3

$ws2.Range("Q37").Value  = 3

# --- Sheet2: roll-up rows (COUNT / SUM per location) -------------------
$ws2.Range("A40").Value = "Total"
$ws2.Range("Q40").Formula    = "=COUNT(Q3:Q39)"
$ws2.Range("R40:AA40").Formula = "=COUNT(R3:R39)"
$ws2.Range("Q41").Formula    = "=SUM(Q3:Q39)"
$ws2.Range("R41:AA41").Formula = "=SUM(R3:R39)"

# --- restore on-screen selections (Sheet1 stays frozen/scrolled as-is,
#     Sheet2 remains the active tab) ------------------------------------
$ws1.Activate()
$ws1.Range("L16").Select()
$ws2.Activate()
$ws2.Range("Q43").Select()
